# "Generate Report for Handoff"
# The report rows that describe the two handback files swap places:
# the file that used to be "a474ad6f..." now reports as "8d52b95f..." and
# vice-versa, and the 8d52b95f row moves from "Handed back" to
# "Ready for handoff" (stale handback) while a474ad6f keeps being in sync.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.md"
$ov.Range("B2").Value = "e2e\a474ad6f-392e-4361-909c-9ce03469b8ae.md"

$ov.Range("A3").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.md"
$ov.Range("B3").Value = "e2e\8d52b95f-5e77-4b29-9701-277e4ab73e11.md"

$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-09-07 05:50:51"

# Rebuild the hyperlinks on column B so the displayed text matches the new
# file names while keeping the same target addresses (same rIds order).
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61d1f749269eb7372d4a855d2d80612a6be1a328/e2e/8d52b95f-5e77-4b29-9701-277e4ab73e11.md", "", "", "e2e\a474ad6f-392e-4361-909c-9ce03469b8ae.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61d1f749269eb7372d4a855d2d80612a6be1a328/e2e/a474ad6f-392e-4361-909c-9ce03469b8ae.md", "", "", "e2e\8d52b95f-5e77-4b29-9701-277e4ab73e11.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.md"
$zh.Range("G2").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.4124c72963f29ffa7603a473b77e015cf4fd79dc.zh-cn.xlf"
$zh.Range("I2").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.md"
$zh.Range("J2").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.4124c72963f29ffa7603a473b77e015cf4fd79dc.zh-cn.xlf"

$zh.Range("A3").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.93c95854e16b1616a0761d7d9acba8bd20766fdf.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-07 05:50:40"
$zh.Range("I3").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.md"
$zh.Range("J3").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.93c95854e16b1616a0761d7d9acba8bd20766fdf.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61d1f749269eb7372d4a855d2d80612a6be1a328/e2e/8d52b95f-5e77-4b29-9701-277e4ab73e11.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01a6e2e3fdfc064db544854c6de8f09ebb959ac4/e2e/8d52b95f-5e77-4b29-9701-277e4ab73e11.md."

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61d1f749269eb7372d4a855d2d80612a6be1a328/e2e/8d52b95f-5e77-4b29-9701-277e4ab73e11.md", "", "", "a474ad6f-392e-4361-909c-9ce03469b8ae.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ec93d39add3e95707bb66ffbf5f7a902df5443a6/e2e/8d52b95f-5e77-4b29-9701-277e4ab73e11.md", "", "", "a474ad6f-392e-4361-909c-9ce03469b8ae.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61d1f749269eb7372d4a855d2d80612a6be1a328/e2e/a474ad6f-392e-4361-909c-9ce03469b8ae.md", "", "", "8d52b95f-5e77-4b29-9701-277e4ab73e11.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ec93d39add3e95707bb66ffbf5f7a902df5443a6/e2e/a474ad6f-392e-4361-909c-9ce03469b8ae.md", "", "", "8d52b95f-5e77-4b29-9701-277e4ab73e11.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.md"
$de.Range("G2").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.4124c72963f29ffa7603a473b77e015cf4fd79dc.de-de.xlf"
$de.Range("I2").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.md"
$de.Range("J2").Value = "a474ad6f-392e-4361-909c-9ce03469b8ae.4124c72963f29ffa7603a473b77e015cf4fd79dc.de-de.xlf"

$de.Range("A3").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.93c95854e16b1616a0761d7d9acba8bd20766fdf.de-de.xlf"
$de.Range("H3").Value = "2016-09-07 05:50:51"
$de.Range("I3").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.md"
$de.Range("J3").Value = "8d52b95f-5e77-4b29-9701-277e4ab73e11.93c95854e16b1616a0761d7d9acba8bd20766fdf.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61d1f749269eb7372d4a855d2d80612a6be1a328/e2e/8d52b95f-5e77-4b29-9701-277e4ab73e11.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01a6e2e3fdfc064db544854c6de8f09ebb959ac4/e2e/8d52b95f-5e77-4b29-9701-277e4ab73e11.md."

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61d1f749269eb7372d4a855d2d80612a6be1a328/e2e/8d52b95f-5e77-4b29-9701-277e4ab73e11.md", "", "", "a474ad6f-392e-4361-909c-9ce03469b8ae.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/da2cb78011401147c193ff0997cdbd487e132c06/e2e/8d52b95f-5e77-4b29-9701-277e4ab73e11.md", "", "", "a474ad6f-392e-4361-909c-9ce03469b8ae.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61d1f749269eb7372d4a855d2d80612a6be1a328/e2e/a474ad6f-392e-4361-909c-9ce03469b8ae.md", "", "", "8d52b95f-5e77-4b29-9701-277e4ab73e11.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/da2cb78011401147c193ff0997cdbd487e132c06/e2e/a474ad6f-392e-4361-909c-9ce03469b8ae.md", "", "", "8d52b95f-5e77-4b29-9701-277e4ab73e11.md") | Out-Null

$wb.Save()
